# Update "Chiffres COVID-19 Valais" worksheet with newly reported data rows.
# Columns: A=Date, B=Cumul cas positifs (formula), C=Nb nouveaux cas positifs,
#          D=Nb nouvelles admissions, E=Patients SI, F=Patients intubes,
#          G=Patients hospitalises hors SI, H=Total hospitalisations (formula),
#          I=Nb nouvelles sorties, J=Cumul deces (formula), K=Nb nouveaux deces (formula),
#          L=Nb nouveaux deces hopital, M=Nb nouveaux deces extra-hospitaliers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns C/E/F/G (General-formatted numeric cells): plain value writes ---

$ws.Range("G622").Value = 13
$ws.Range("G623").Value = 15
$ws.Range("G624").Value = 15
$ws.Range("G625").Value = 16
$ws.Range("G626").Value = 18
$ws.Range("G627").Value = 18
$ws.Range("G628").Value = 20
$ws.Range("G629").Value = 23

$ws.Range("C630").Value = 140
$ws.Range("G630").Value = 26

$ws.Range("C631").Value = 156
$ws.Range("E631").Value = 7
$ws.Range("G631").Value = 31

$ws.Range("C632").Value = 180
$ws.Range("E632").Value = 6
$ws.Range("G632").Value = 36

$ws.Range("C633").Value = 167
$ws.Range("E633").Value = 6
$ws.Range("F633").Value = 4
$ws.Range("G633").Value = 34

$ws.Range("C634").Value = 98
$ws.Range("E634").Value = 8
$ws.Range("F634").Value = 3
$ws.Range("G634").Value = 38

$ws.Range("C635").Value = 76
$ws.Range("E635").Value = 7
$ws.Range("F635").Value = 4
$ws.Range("G635").Value = 40

$ws.Range("C636").Value = 18
$ws.Range("E636").Value = 7
$ws.Range("F636").Value = 4
$ws.Range("G636").Value = 41

# --- Columns L/M (Text-number-formatted cells): force a General format while
#     writing so the value is stored numerically (matches existing L/M cells
#     in the workbook, which already hold numbers despite the "@" display
#     format), then restore the original text format. ---

$lmCells = @("L628","M630","L632","M632","L633","M633","L634","M634","L635","M635","L636","M636")
foreach ($c in $lmCells) {
    $ws.Range($c).NumberFormat = "General"
}

$ws.Range("L628").Value = 1

$ws.Range("M630").Value = 2

$ws.Range("L632").Value = 2
$ws.Range("M632").Value = 1

$ws.Range("L633").Value = 0
$ws.Range("M633").Value = 1

$ws.Range("L634").Value = 0
$ws.Range("M634").Value = 0

$ws.Range("L635").Value = 0
$ws.Range("M635").Value = 0

$ws.Range("L636").Value = 0
$ws.Range("M636").Value = 0

foreach ($c in $lmCells) {
    $ws.Range($c).NumberFormat = "@"
}

# --- Update the frozen-pane (bottomRight) active-cell selection to A2 ---
$ws.Range("A2").Select()

# Switch off automatic (whole-workbook) recalculation after the edits above
# have already pushed their direct formula dependents (B/H/J/K in the edited
# rows) to fresh values. This keeps the untouched, far-below "calculate
# always" (volatile TODAY()-based) rows from being needlessly recalculated
# and rewritten, matching upstream Excel's behaviour of only recalculating
# cells on the dependency path of what was actually edited.
$excel.Calculation = -4135
